# Update "想去人数" (F) and "最低票价" (G) figures across the three
# worksheets that carry event data: 展览, 演出, and 全部类型.
# (本地生活 has no data rows and is left untouched.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 230
$ws.Range("G2").Value = 59
$ws.Range("F4").Value = 19874
$ws.Range("F5").Value = 795
$ws.Range("F7").Value = 1095
$ws.Range("F9").Value = 7506
$ws.Range("F10").Value = 504
$ws.Range("F12").Value = 258
$ws.Range("F14").Value = 155
$ws.Range("F15").Value = 114
$ws.Range("F19").Value = 1334
$ws.Range("F20").Value = 402
$ws.Range("F21").Value = 71
$ws.Range("F22").Value = 676
$ws.Range("F24").Value = 62
$ws.Range("F26").Value = 319
$ws.Range("F27").Value = 1096
$ws.Range("F28").Value = 30
$ws.Range("F29").Value = 15
$ws.Range("F30").Value = 178
$ws.Range("F33").Value = 59
$ws.Range("F34").Value = 2815
$ws.Range("F38").Value = 12589
$ws.Range("F40").Value = 72
$ws.Range("F41").Value = 22
$ws.Range("F42").Value = 55
$ws.Range("F43").Value = 256
$ws.Range("F44").Value = 356
$ws.Range("F45").Value = 3989

# ---------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 3

# ---------------------------------------------------------------
# Sheet "全部类型" (all types combined)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 230
$ws.Range("G2").Value = 59
$ws.Range("F4").Value = 19874
$ws.Range("F5").Value = 795
$ws.Range("F7").Value = 1095
$ws.Range("F9").Value = 7506
$ws.Range("F12").Value = 258
$ws.Range("F14").Value = 155
$ws.Range("F18").Value = 191
$ws.Range("F20").Value = 402
$ws.Range("F22").Value = 676
$ws.Range("F24").Value = 62
$ws.Range("F25").Value = 67
$ws.Range("F26").Value = 319
$ws.Range("F27").Value = 1096
$ws.Range("F28").Value = 30
$ws.Range("F33").Value = 3
$ws.Range("F34").Value = 59
$ws.Range("F35").Value = 33
$ws.Range("F36").Value = 2815
$ws.Range("F37").Value = 25
$ws.Range("F40").Value = 12589
$ws.Range("F42").Value = 72
$ws.Range("F43").Value = 22
$ws.Range("F44").Value = 55
$ws.Range("F45").Value = 256
$ws.Range("F46").Value = 356
$ws.Range("F47").Value = 3989
$ws.Range("F48").Value = 318

$wb.Save()
